$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Cells.Item(9, 1).Value = 112171801
$ws.Cells.Item(9, 2).Value = 78699
$ws.Cells.Item(9, 4).Value = 'NT'
$ws.Cells.Item(9, 5).Value = 6458
$ws.Cells.Item(9, 6).Value = 'Lunglav'
$ws.Cells.Item(9, 7).Value = 'Lobaria pulmonaria'
$ws.Cells.Item(9, 8).Value = '(L.) Hoffm.'
$ws.Cells.Item(9, 17).Value = 756448
$ws.Cells.Item(9, 18).Value = 7212052
$ws.Cells.Item(9, 36).ClearContents()
$ws.Cells.Item(9, 37).ClearContents()
$ws.Cells.Item(9, 41).ClearContents()

# Row 10
$ws.Cells.Item(10, 1).Value = 112171776
$ws.Cells.Item(10, 2).Value = 85434
$ws.Cells.Item(10, 4).Value = 'NT'
$ws.Cells.Item(10, 5).Value = 3739
$ws.Cells.Item(10, 6).Value = 'Persiljespindling'
$ws.Cells.Item(10, 7).Value = 'Cortinarius sulfurinus'
$ws.Cells.Item(10, 8).Value = 'Quél.'
$ws.Cells.Item(10, 17).Value = 756261
$ws.Cells.Item(10, 18).Value = 7211953

# Row 11
$ws.Cells.Item(11, 1).Value = 112171779
$ws.Cells.Item(11, 2).Value = 102166
$ws.Cells.Item(11, 4).Value = 'LC'
$ws.Cells.Item(11, 5).Value = 222412
$ws.Cells.Item(11, 6).Value = 'Tibast'
$ws.Cells.Item(11, 7).Value = 'Daphne mezereum'
$ws.Cells.Item(11, 8).Value = 'L.'
$ws.Cells.Item(11, 17).Value = 756291
$ws.Cells.Item(11, 18).Value = 7211892

# Row 12
$ws.Cells.Item(12, 1).Value = 112171814
$ws.Cells.Item(12, 2).Value = 89485
$ws.Cells.Item(12, 5).Value = 112
$ws.Cells.Item(12, 6).Value = 'Stjärntagging'
$ws.Cells.Item(12, 7).Value = 'Asterodon ferruginosus'
$ws.Cells.Item(12, 8).Value = 'Pat.'
$ws.Cells.Item(12, 17).Value = 756486
$ws.Cells.Item(12, 18).Value = 7212041
$ws.Cells.Item(12, 36).Value = 'gran'
$ws.Cells.Item(12, 37).Value = 'Picea abies'
$ws.Cells.Item(12, 41).Value = 'Picea abies'

# Row 13
$ws.Cells.Item(13, 1).Value = 112171812
$ws.Cells.Item(13, 2).Value = 78732
$ws.Cells.Item(13, 4).Value = 'LC'
$ws.Cells.Item(13, 5).Value = 6463
$ws.Cells.Item(13, 6).Value = 'Bårdlav'
$ws.Cells.Item(13, 7).Value = 'Nephroma parile'
$ws.Cells.Item(13, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(13, 17).Value = 756485
$ws.Cells.Item(13, 18).Value = 7212023
$ws.Cells.Item(13, 36).Value = 'sälg'
$ws.Cells.Item(13, 37).Value = 'Salix caprea'
$ws.Cells.Item(13, 41).Value = 'Salix caprea'

# Row 14
$ws.Cells.Item(14, 1).Value = 112171806
$ws.Cells.Item(14, 2).Value = 86357
$ws.Cells.Item(14, 5).Value = 4412
$ws.Cells.Item(14, 6).Value = 'Äggvaxskivling'
$ws.Cells.Item(14, 7).Value = 'Hygrophorus karstenii'
$ws.Cells.Item(14, 8).Value = 'Sacc. & Cub.'
$ws.Cells.Item(14, 17).Value = 756477
$ws.Cells.Item(14, 18).Value = 7212031

# Row 15
$ws.Cells.Item(15, 1).Value = 112171795
$ws.Cells.Item(15, 2).Value = 77636
$ws.Cells.Item(15, 5).Value = 6425
$ws.Cells.Item(15, 6).Value = 'Garnlav'
$ws.Cells.Item(15, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(15, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(15, 17).Value = 756378
$ws.Cells.Item(15, 18).Value = 7212050

# Row 16
$ws.Cells.Item(16, 1).Value = 112171785
$ws.Cells.Item(16, 2).Value = 78663
$ws.Cells.Item(16, 4).Value = 'LC'
$ws.Cells.Item(16, 5).Value = 229748
$ws.Cells.Item(16, 6).Value = 'Gytterlav'
$ws.Cells.Item(16, 7).Value = 'Protopannaria pezizoides'
$ws.Cells.Item(16, 8).Value = '(Weber) P.M.Jørg. & S.Ekman'
$ws.Cells.Item(16, 17).Value = 756412
$ws.Cells.Item(16, 29).Value = 'på berg'
$ws.Cells.Item(16, 36).ClearContents()
$ws.Cells.Item(16, 37).ClearContents()
$ws.Cells.Item(16, 41).ClearContents()

# Row 17
$ws.Cells.Item(17, 1).Value = 112171792
$ws.Cells.Item(17, 2).Value = 85387
$ws.Cells.Item(17, 5).Value = 249228
$ws.Cells.Item(17, 6).Value = 'Barrfagerspindling'
$ws.Cells.Item(17, 7).Value = 'Cortinarius piceae'
$ws.Cells.Item(17, 8).Value = 'Frøslev, T.S.Jeppesen & Brandrud'
$ws.Cells.Item(17, 17).Value = 756395
$ws.Cells.Item(17, 18).Value = 7211974
$ws.Cells.Item(17, 29).ClearContents()

# Row 18
$ws.Cells.Item(18, 1).Value = 112171798
$ws.Cells.Item(18, 2).Value = 78699
$ws.Cells.Item(18, 17).Value = 756371
$ws.Cells.Item(18, 18).Value = 7212116
$ws.Cells.Item(18, 36).ClearContents()
$ws.Cells.Item(18, 37).ClearContents()
$ws.Cells.Item(18, 41).ClearContents()

# Row 19
$ws.Cells.Item(19, 1).Value = 112171813
$ws.Cells.Item(19, 2).Value = 78699
$ws.Cells.Item(19, 5).Value = 6458
$ws.Cells.Item(19, 6).Value = 'Lunglav'
$ws.Cells.Item(19, 7).Value = 'Lobaria pulmonaria'
$ws.Cells.Item(19, 8).Value = '(L.) Hoffm.'
$ws.Cells.Item(19, 17).Value = 756485
$ws.Cells.Item(19, 18).Value = 7212023
$ws.Cells.Item(19, 36).Value = 'sälg'
$ws.Cells.Item(19, 37).Value = 'Salix caprea'
$ws.Cells.Item(19, 41).Value = 'Salix caprea'

# Row 20
$ws.Cells.Item(20, 1).Value = 112171788
$ws.Cells.Item(20, 2).Value = 78699
$ws.Cells.Item(20, 5).Value = 6458
$ws.Cells.Item(20, 6).Value = 'Lunglav'
$ws.Cells.Item(20, 7).Value = 'Lobaria pulmonaria'
$ws.Cells.Item(20, 8).Value = '(L.) Hoffm.'
$ws.Cells.Item(20, 17).Value = 756401
$ws.Cells.Item(20, 18).Value = 7211954
$ws.Cells.Item(20, 36).Value = 'sälg'
$ws.Cells.Item(20, 37).Value = 'Salix caprea'
$ws.Cells.Item(20, 41).Value = 'Salix caprea'

# Row 21
$ws.Cells.Item(21, 1).Value = 112171787
$ws.Cells.Item(21, 2).Value = 86357
$ws.Cells.Item(21, 4).Value = 'NT'
$ws.Cells.Item(21, 5).Value = 4412
$ws.Cells.Item(21, 6).Value = 'Äggvaxskivling'
$ws.Cells.Item(21, 7).Value = 'Hygrophorus karstenii'
$ws.Cells.Item(21, 8).Value = 'Sacc. & Cub.'
$ws.Cells.Item(21, 17).Value = 756408
$ws.Cells.Item(21, 18).Value = 7211956

# Row 22
$ws.Cells.Item(22, 1).Value = 112171810
$ws.Cells.Item(22, 2).Value = 90466
$ws.Cells.Item(22, 5).Value = 4769
$ws.Cells.Item(22, 6).Value = 'Svavelriska'
$ws.Cells.Item(22, 7).Value = 'Lactarius scrobiculatus'
$ws.Cells.Item(22, 8).Value = '(Scop.:Fr.) Fr.'
$ws.Cells.Item(22, 17).Value = 756486
$ws.Cells.Item(22, 18).Value = 7212020
